$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update production_date values in column D (rows 2-6) to new relative dates.
# Force text number format so Excel does not reinterpret the strings as date serials.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2026-02-12"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2026-02-13"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2026-02-14"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2026-02-15"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2026-02-16"
